# Auto-generated Excel COM-interop script to apply diff changes
# "Update latest output (run 166)"

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("B2").Value = 46062.25
$wsSchedule.Range("C2").Value = 6
$wsSchedule.Range("D2").Value = 22.68
$wsSchedule.Range("E2").Value = 1244.68091475
$wsSchedule.Range("F2").Value = 54.88011087962964
$wsSchedule.Range("A3").Value = 46062.33333333334
$wsSchedule.Range("C3").Value = 8
$wsSchedule.Range("D3").Value = 30.24
$wsSchedule.Range("E3").Value = 1439.45284275
$wsSchedule.Range("F3").Value = 47.60095379464285
$wsSchedule.Range("A4").Value = 46062.91666666666
$wsSchedule.Range("B4").Value = 46063.16666666666
$wsSchedule.Range("C4").Value = 6
$wsSchedule.Range("D4").Value = 22.68
$wsSchedule.Range("E4").Value = 1251.0008745
$wsSchedule.Range("F4").Value = 55.15876871693122
$wsSchedule.Range("B5").Value = 46063.64583333334
$wsSchedule.Range("C5").Value = 8
$wsSchedule.Range("D5").Value = 30.24
$wsSchedule.Range("E5").Value = 976.4739659999998
$wsSchedule.Range("F5").Value = 32.29080575396825

$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("E13").Value = "ON"
$wsDetailed.Range("E18").Value = "ON"
$wsDetailed.Range("B38").Value = 115
$wsDetailed.Range("B39").Value = 130.69454
$wsDetailed.Range("B40").Value = 166.99
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 216.91053
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 218.70531
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 212.6158
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 166.99
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("E44").Value = "OFF"
$wsDetailed.Range("B45").Value = 147.51
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 116.00628
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 138.42
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 108.01
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 105
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B50").Value = 105.79
$wsDetailed.Range("B51").Value = 105.68919
$wsDetailed.Range("B52").Value = 103.32057
$wsDetailed.Range("B53").Value = 103.30972
$wsDetailed.Range("B54").Value = 105.63067
$wsDetailed.Range("B55").Value = 102.8976
$wsDetailed.Range("B56").Value = 102.9661
$wsDetailed.Range("B57").Value = 86.03769
$wsDetailed.Range("E57").Value = "ON"
$wsDetailed.Range("B58").Value = 101.86026
$wsDetailed.Range("B59").Value = 123.44747
$wsDetailed.Range("B60").Value = 121.10472
$wsDetailed.Range("B61").Value = 108.01
$wsDetailed.Range("B62").Value = 111.90666
$wsDetailed.Range("B63").Value = 113.32199
$wsDetailed.Range("B66").Value = 76.94002
$wsDetailed.Range("B67").Value = 64.44902999999999
$wsDetailed.Range("B69").Value = 36.07
$wsDetailed.Range("B70").Value = 47.06354
$wsDetailed.Range("B71").Value = 47.04111
$wsDetailed.Range("B72").Value = 47.61559
$wsDetailed.Range("B73").Value = 55.63646
$wsDetailed.Range("B74").Value = 55.30004
$wsDetailed.Range("B75").Value = 62.03642
$wsDetailed.Range("B76").Value = 61.57854
$wsDetailed.Range("B77").Value = 64.88755
$wsDetailed.Range("B78").Value = 69.90346
$wsDetailed.Range("B80").Value = 100.01
$wsDetailed.Range("B81").Value = 88.81128
$wsDetailed.Range("E81").Value = "OFF"
$wsDetailed.Range("B82").Value = 69.98341000000001
$wsDetailed.Range("B83").Value = 69.92310000000001
$wsDetailed.Range("B84").Value = 70.16486
$wsDetailed.Range("B85").Value = 108.45901
$wsDetailed.Range("B86").Value = 110.17505
$wsDetailed.Range("B87").Value = 161.25495
$wsDetailed.Range("B88").Value = 206.75779
$wsDetailed.Range("B91").Value = 222.23033
$wsDetailed.Range("B92").Value = 178.67423
$wsDetailed.Range("B93").Value = 166.99
$wsDetailed.Range("B94").Value = 144.62829
$wsDetailed.Range("B95").Value = 154.2
$wsDetailed.Range("B96").Value = 154.62606
$wsDetailed.Range("B97").Value = 138.42
